# Scheduled market-board price refresh for the Leve profit tracker.
# Updates currentAveragePrice(NQ/HQ), LevePrice(NQ/HQ) and the derived
# LeveProfit(NQ/HQ) columns (H-N) for the leves whose item prices moved
# since the last run. One worksheet per Disciple of the Hand job.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 8: On the Drip (Eye Drops)
$ws.Range("H8").Value = 12517.75
$ws.Range("I8").Value = 23.666666
$ws.Range("K8").Value = 70.99999800000001
$ws.Range("M8").Value = 68.00000199999999

# Row 34: Sophomore Slump (Goatskin Grimoire)
$ws.Range("H34").Value = 20655.092
$ws.Range("I34").Value = 20655.092
$ws.Range("K34").Value = 20655.092
$ws.Range("M34").Value = -20452.092

# Row 36: You Put Your Left Hand In (Engraved Goatskin Grimoire)
$ws.Range("H36").Value = 20655.092
$ws.Range("I36").Value = 20655.092
$ws.Range("K36").Value = 20655.092
$ws.Range("M36").Value = -19940.092

# Row 107: Another Man's Ink (Enchanted Truegold Ink)
$ws.Range("H107").Value = 510.43478
$ws.Range("I107").Value = 533.6842
$ws.Range("K107").Value = 533.6842
$ws.Range("M107").Value = 1386.3158

# Row 129: Practical Command (Commanding Craftsman's Draught)
$ws.Range("H129").Value = 980.8261
$ws.Range("I129").Value = 526.9
$ws.Range("J129").Value = 1057.7627
$ws.Range("K129").Value = 1580.7
$ws.Range("L129").Value = 3173.2881
$ws.Range("M129").Value = 3419.3
$ws.Range("N129").Value = -13173.2881

# Row 132: Fast-forwarding Flora (Growth Formula Lambda)
$ws.Range("H132").Value = 1887.3889
$ws.Range("I132").Value = 1612.2709
$ws.Range("K132").Value = 4836.8127
$ws.Range("M132").Value = -2306.8127

# Row 137: Cutting Edge of Culinary Quality (Magnesia Whetstone)
$ws.Range("H137").Value = 1521.0968
$ws.Range("I137").Value = 1315.6222
$ws.Range("J137").Value = 2065
$ws.Range("K137").Value = 3946.8666
$ws.Range("L137").Value = 6195
$ws.Range("M137").Value = -1396.8666
$ws.Range("N137").Value = -11295

$ws = $wb.Worksheets.Item("ARM")
# Row 2: Ain't Got No Ingots (Bronze Ingot)
$ws.Range("H2").Value = 859.2162
$ws.Range("I2").Value = 579.4666999999999
$ws.Range("J2").Value = 2058.1428
$ws.Range("K2").Value = 579.4666999999999
$ws.Range("L2").Value = 2058.1428
$ws.Range("M2").Value = -466.4666999999999
$ws.Range("N2").Value = -2284.1428

# Row 61: Dealing with the Tough Stuff (Cobalt Ingot)
$ws.Range("H61").Value = 1173.9556
$ws.Range("I61").Value = 1112.3103
$ws.Range("J61").Value = 1285.6875
$ws.Range("K61").Value = 1112.3103
$ws.Range("L61").Value = 1285.6875
$ws.Range("M61").Value = -900.3103000000001
$ws.Range("N61").Value = -1709.6875

# Row 63: Rivets Run through It (Mythrite Rivets)
$ws.Range("H63").Value = 6346.357
$ws.Range("I63").Value = 5194.4443
$ws.Range("J63").Value = 8419.799999999999
$ws.Range("K63").Value = 5194.4443
$ws.Range("L63").Value = 8419.799999999999
$ws.Range("M63").Value = -4508.4443
$ws.Range("N63").Value = -9791.799999999999

# Row 66: A Riveting Revival (L) (Mythrite Rivets)
$ws.Range("H66").Value = 6346.357
$ws.Range("I66").Value = 5194.4443
$ws.Range("J66").Value = 8419.799999999999
$ws.Range("K66").Value = 25972.2215
$ws.Range("L66").Value = 42099
$ws.Range("M66").Value = -22540.2215
$ws.Range("N66").Value = -48963

# Row 116: No Scope (Titanbronze Ingot)
$ws.Range("H116").Value = 859.2162
$ws.Range("I116").Value = 579.4666999999999
$ws.Range("J116").Value = 2058.1428
$ws.Range("K116").Value = 579.4666999999999
$ws.Range("L116").Value = 2058.1428
$ws.Range("M116").Value = 1714.5333
$ws.Range("N116").Value = -6646.1428

# Row 136: Metal with Mettle (Cobalt Tungsten Ingot)
$ws.Range("H136").Value = 1173.9556
$ws.Range("I136").Value = 1112.3103
$ws.Range("J136").Value = 1285.6875
$ws.Range("K136").Value = 3336.9309
$ws.Range("L136").Value = 3857.0625
$ws.Range("M136").Value = -786.9309000000003
$ws.Range("N136").Value = -8957.0625

$ws = $wb.Worksheets.Item("BSM")
# Row 3: Hells Bells (Bronze Ingot)
$ws.Range("H3").Value = 859.2162
$ws.Range("I3").Value = 579.4666999999999
$ws.Range("J3").Value = 2058.1428
$ws.Range("K3").Value = 579.4666999999999
$ws.Range("L3").Value = 2058.1428
$ws.Range("M3").Value = -465.4666999999999
$ws.Range("N3").Value = -2286.1428

# Row 15: Anutha Spatha (Bronze Spatha)
$ws.Range("H15").Value = 23777.777
$ws.Range("I15").Value = 23428.572
$ws.Range("J15").Value = 25000
$ws.Range("K15").Value = 23428.572
$ws.Range("L15").Value = 25000
$ws.Range("M15").Value = -23201.572
$ws.Range("N15").Value = -25454

# Row 82: Spirituality Inspector (Titanium Lump Hammer)
$ws.Range("H82").Value = 58897
$ws.Range("I82").Value = 59886.7
$ws.Range("K82").Value = 59886.7
$ws.Range("M82").Value = -59503.7

# Row 85: The Clamor for Hammers (L) (Titanium Lump Hammer)
$ws.Range("H85").Value = 58897
$ws.Range("I85").Value = 59886.7
$ws.Range("K85").Value = 59886.7
$ws.Range("M85").Value = -58560.7

$ws = $wb.Worksheets.Item("CRP")
# Row 15: On the Move (Ragstone Grinding Wheel)
$ws.Range("H15").Value = 24377.25
$ws.Range("J15").Value = 24377.25
$ws.Range("L15").Value = 24377.25
$ws.Range("N15").Value = -24717.25

# Row 58: You Do the Heavy Lifting (Mahogany Lumber)
$ws.Range("H58").Value = 1563.3077
$ws.Range("I58").Value = 1524.8182
$ws.Range("J58").Value = 1775
$ws.Range("K58").Value = 1524.8182
$ws.Range("L58").Value = 1775
$ws.Range("M58").Value = -1321.8182
$ws.Range("N58").Value = -2181

# Row 120: Kindling the Flame (Lignum Vitae Ring)
$ws.Range("H120").Value = 65976
$ws.Range("J120").Value = 65976
$ws.Range("L120").Value = 65976
$ws.Range("N120").Value = -73234

# Row 132: Hull Lotta Damage (Ginseng Lumber)
$ws.Range("H132").Value = 2265.0488
$ws.Range("I132").Value = 2077.162
$ws.Range("K132").Value = 6231.485999999999
$ws.Range("M132").Value = -3701.485999999999

# Row 134: Wood You Be Quiet (Ceiba Lumber)
$ws.Range("H134").Value = 1586.6562
$ws.Range("I134").Value = 1412.9546
$ws.Range("J134").Value = 1968.8
$ws.Range("K134").Value = 4238.8638
$ws.Range("L134").Value = 5906.4
$ws.Range("M134").Value = -1703.8638
$ws.Range("N134").Value = -10976.4

# Row 136: Turali Quality (Dark Mahogany Lumber)
$ws.Range("H136").Value = 1563.3077
$ws.Range("I136").Value = 1524.8182
$ws.Range("J136").Value = 1775
$ws.Range("K136").Value = 4574.4546
$ws.Range("L136").Value = 5325
$ws.Range("M136").Value = -2024.4546
$ws.Range("N136").Value = -10425

$ws = $wb.Worksheets.Item("CUL")
# Row 3: Trout Fishing in Limsa (Grilled Trout)
$ws.Range("H3").Value = 4898.5713
$ws.Range("I3").Value = 4048.3333
$ws.Range("K3").Value = 12144.9999
$ws.Range("M3").Value = -12032.9999

# Row 132: More Mezcal (Cooking Mezcal)
$ws.Range("H132").Value = 2082.3928
$ws.Range("J132").Value = 2653.2942
$ws.Range("L132").Value = 23879.6478
$ws.Range("N132").Value = -28939.6478

$ws = $wb.Worksheets.Item("GSM")
# Row 109: You're My Wonderhall (Hematite Earrings of Healing)
$ws.Range("H109").Value = 20659.25
$ws.Range("J109").Value = 20659.25
$ws.Range("L109").Value = 20659.25
$ws.Range("N109").Value = -22739.25

# Row 113: Copious Crystal Cannons (Manasilver Nugget)
$ws.Range("H113").Value = 1598.6154
$ws.Range("I113").Value = 898.6667
$ws.Range("J113").Value = 2198.5715
$ws.Range("K113").Value = 898.6667
$ws.Range("L113").Value = 2198.5715
$ws.Range("M113").Value = 1271.3333
$ws.Range("N113").Value = -6538.5715

# Row 123: Workplace Workout (Ametrine Ring of Fending)
$ws.Range("H123").Value = 15755.409
$ws.Range("J123").Value = 15755.409
$ws.Range("L123").Value = 15755.409
$ws.Range("N123").Value = -20655.409

# Row 132: On Board for Lar (Lar Ingot)
$ws.Range("H132").Value = 2878.8333
$ws.Range("I132").Value = 2514.8235
$ws.Range("J132").Value = 3762.8572
$ws.Range("K132").Value = 7544.470499999999
$ws.Range("L132").Value = 11288.5716
$ws.Range("M132").Value = -5014.470499999999
$ws.Range("N132").Value = -16348.5716

$ws = $wb.Worksheets.Item("LTW")
# Row 22: Skin off Their Backs (Aldgoat Leather)
$ws.Range("H22").Value = 1072.7273
$ws.Range("I22").Value = 2120
$ws.Range("J22").Value = 200
$ws.Range("K22").Value = 2120
$ws.Range("L22").Value = 200
$ws.Range("M22").Value = -1825
$ws.Range("N22").Value = -790

# Row 27: Fire and Hide (Aldgoat Leather)
$ws.Range("H27").Value = 1072.7273
$ws.Range("I27").Value = 2120
$ws.Range("J27").Value = 200
$ws.Range("K27").Value = 2120
$ws.Range("L27").Value = 200
$ws.Range("M27").Value = -2013
$ws.Range("N27").Value = -414

# Row 136: Respect for Br'aax (Br'aax Leather)
$ws.Range("H136").Value = 3424.52
$ws.Range("I136").Value = 3636.3022
$ws.Range("J136").Value = 2123.5715
$ws.Range("K136").Value = 10908.9066
$ws.Range("L136").Value = 6370.7145
$ws.Range("M136").Value = -8358.9066
$ws.Range("N136").Value = -11470.7145

$ws = $wb.Worksheets.Item("WVR")
# Row 132: Comfy Cabins (Snow Cotton Cloth)
$ws.Range("H132").Value = 2521.5652
$ws.Range("I132").Value = 1764.1765
$ws.Range("J132").Value = 4667.5
$ws.Range("K132").Value = 5292.529500000001
$ws.Range("L132").Value = 14002.5
$ws.Range("M132").Value = -2762.529500000001
$ws.Range("N132").Value = -19062.5

# Row 136: Weaving the Envelope (Sarcenet Cloth)
$ws.Range("H136").Value = 1274.8591
$ws.Range("I136").Value = 1242.3396
$ws.Range("J136").Value = 1370.6111
$ws.Range("K136").Value = 3727.0188
$ws.Range("L136").Value = 4111.8333
$ws.Range("M136").Value = -1177.0188
$ws.Range("N136").Value = -9211.8333
